$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fixes: row 1 ---
$ws.Range("Q1").Value = "Kevin Lee"

# New columns AE:AH in header row, copying style from AD1 (bold, centered, bordered)
$ws.Range("AD1").Copy($ws.Range("AE1"))
$ws.Range("AD1").Copy($ws.Range("AF1"))
$ws.Range("AD1").Copy($ws.Range("AG1"))
$ws.Range("AD1").Copy($ws.Range("AH1"))
$ws.Range("AE1").Value = "Spencer Harris"
$ws.Range("AF1").Value = "Bronte Sundstrom"
$ws.Range("AG1").Value = "Parker Simpson"
$ws.Range("AH1").Value = "Abby LeahFriend"

# --- New rows 104:115 ("beach week" and "the last roll" games) ---

# Row 104
$ws.Range("A103").Copy($ws.Range("A104"))
$ws.Range("A104").Value = 103
$ws.Range("F104").Value = 1544.302477264692
$ws.Range("H104").Value = 1206.063981098702
$ws.Range("I104").Value = 979.4502013604686
$ws.Range("AF104").Value = 1192.819513605344

# Row 105
$ws.Range("A104").Copy($ws.Range("A105"))
$ws.Range("A105").Value = 104
$ws.Range("F105").Value = 1545.742079244441
$ws.Range("G105").Value = 1269.209554187589
$ws.Range("H105").Value = 1207.503583078451
$ws.Range("M105").Value = 1137.252855431793

# Row 106
$ws.Range("A105").Copy($ws.Range("A106"))
$ws.Range("A106").Value = 105
$ws.Range("D106").Value = 1313.684825654176
$ws.Range("Q106").Value = 1427.419378002535
$ws.Range("T106").Value = 1162.43448193746
$ws.Range("Z106").Value = 1040.370685155026

# Row 107
$ws.Range("A106").Copy($ws.Range("A107"))
$ws.Range("A107").Value = 106
$ws.Range("F107").Value = 1546.311365205288
$ws.Range("H107").Value = 1208.072869039298
$ws.Range("R107").Value = 1183.082625265769
$ws.Range("AE107").Value = 1199.430714039153

# Row 108
$ws.Range("A107").Copy($ws.Range("A108"))
$ws.Range("A108").Value = 107
$ws.Range("D108").Value = 1338.099046793333
$ws.Range("F108").Value = 1521.897144066131
$ws.Range("H108").Value = 1183.658647900141
$ws.Range("Q108").Value = 1451.833599141692

# Row 109
$ws.Range("A108").Copy($ws.Range("A109"))
$ws.Range("A109").Value = 108
$ws.Range("F109").Value = 1498.287424819952
$ws.Range("H109").Value = 1160.048928653962
$ws.Range("M109").Value = 1160.862574677972
$ws.Range("AF109").Value = 1216.429232851523

# Row 110
$ws.Range("A109").Copy($ws.Range("A110"))
$ws.Range("A110").Value = 109
$ws.Range("D110").Value = 1303.641964193648
$ws.Range("G110").Value = 1303.666636787274
$ws.Range("Q110").Value = 1417.376516542007
$ws.Range("R110").Value = 1217.539707865454

# Row 111
$ws.Range("A110").Copy($ws.Range("A111"))
$ws.Range("A111").Value = 110
$ws.Range("G111").Value = 1275.183515985441
$ws.Range("M111").Value = 1189.345695479805
$ws.Range("R111").Value = 1189.056587063621
$ws.Range("AF111").Value = 1244.912353653355

# Row 112
$ws.Range("A111").Copy($ws.Range("A112"))
$ws.Range("A112").Value = 111
$ws.Range("F112").Value = 1504.772979793352
$ws.Range("G112").Value = 1281.669070958841
$ws.Range("S112").Value = 1389.845978659066
$ws.Range("AG112").Value = 1193.5144450266

# Row 113
$ws.Range("A112").Copy($ws.Range("A113"))
$ws.Range("A113").Value = 112
$ws.Range("F113").Value = 1536.548325813933
$ws.Range("G113").Value = 1313.444416979422
$ws.Range("N113").Value = 1201.692490161548
$ws.Range("Q113").Value = 1385.601170521426

# Row 114
$ws.Range("A113").Copy($ws.Range("A114"))
$ws.Range("A114").Value = 113
$ws.Range("F114").Value = 1558.36774250098
$ws.Range("G114").Value = 1335.263833666469
$ws.Range("M114").Value = 1167.526278792758
$ws.Range("AH114").Value = 1178.180583312953

# Row 115
$ws.Range("A114").Copy($ws.Range("A115"))
$ws.Range("A115").Value = 114
$ws.Range("F115").Value = 1578.068408583175
$ws.Range("G115").Value = 1354.964499748664
$ws.Range("T115").Value = 1142.733815855265
$ws.Range("AE115").Value = 1179.730047956958
